# Generate Report for Handback
#
# For the two real localized files (rows 2 & 3) on each language sheet
# (zh-cn / de-de), record that the handback round-tripped successfully:
#   - Status goes from "Ready for handoff" -> "Handed back: in sync with en-US"
#   - "Latest Target File"   (col E) gets the source .md file (same link as col A)
#   - "Latest Handback File" (col F) gets the handed-back .xlf file (same link as col C)
#   - "Latest Handback DateTime" (col G) gets the handback timestamp
#
# The Overview sheet mirrors the same "Status" text for the two files, since
# it shares the same underlying string.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: mirror the status text change (same cells/columns as
# the per-language sheets, just summarized per source file).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper data: per language sheet, per handed-off file row, the
# handback timestamp to stamp into column G.
# ---------------------------------------------------------------------
$sheetInfo = @(
    @{ Index = 2; Timestamp = "2016-02-22 05:18:48" },
    @{ Index = 3; Timestamp = "2016-02-22 05:19:11" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Index)

    # Snapshot existing hyperlink URLs keyed by the cell they decorate,
    # so the new "Latest Target File" / "Latest Handback File" links can
    # reuse the very same targets as the existing "Source File Name" /
    # "Latest Handoff File" links.
    $urlByCell = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $urlByCell[$hl.Range.Address()] = $hl.Address
    }

    $rows = @(2, 3)
    foreach ($r in $rows) {
        $aAddr = "`$A`$$r"
        $cAddr = "`$C`$$r"

        $aDisplay = $ws.Range("A$r").Text
        $cDisplay = $ws.Range("C$r").Text
        $aUrl = $urlByCell[$aAddr]
        $cUrl = $urlByCell[$cAddr]

        # Status
        $ws.Range("B$r").Value = $newStatus

        # Latest Target File (E) -- same file/link as "Source File Name" (A)
        $ws.Range("E$r").Value = $aDisplay
        $ws.Range("E$r").Style = "Hyperlink"
        if ($aUrl) {
            $ws.Hyperlinks.Add($ws.Range("E$r"), $aUrl, [Type]::Missing, [Type]::Missing, $aDisplay)
        }

        # Latest Handback File (F) -- same file/link as "Latest Handoff File" (C)
        $ws.Range("F$r").Value = $cDisplay
        $ws.Range("F$r").Style = "Hyperlink"
        if ($cUrl) {
            $ws.Hyperlinks.Add($ws.Range("F$r"), $cUrl, [Type]::Missing, [Type]::Missing, $cDisplay)
        }

        # Latest Handback DateTime (G)
        $ws.Range("G$r").Value = $info.Timestamp
    }
}
